# Add data for 2022-09-01:
# - Rename sheet "Through 2022-08-23" -> "Through 2022-08-24"
# - Update header label "2022 (through 08-23)" -> "2022 (through 08-24)"
# - Update August count (I9) 130 -> 132
# - Update Total count (I14) 1101 -> 1103

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-08-24"

$ws.Range("I1").Value = "2022 (through 08-24)"
$ws.Range("I9").Value = 132
$ws.Range("I14").Value = 1103
